$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("max"/1) entirely. This shifts the old column D
# ("prediction"/o__Ozemobacterales) into C and the old column E
# ("rejection-f"/o__Ozemobacterales) into D.
$ws.Columns.Item(3).Delete()

# Update the numeric value in B2 (now the only remaining change needed).
$ws.Range("B2").Value = 32637401.27964602
